$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Remove the data row for matricula GO338044453 (ANTONIO WANDERSON DA CONCEIÇÃO),
# shifting all following rows up by one.
$ws.Rows.Item(564).Delete()
